$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price and Volume columns retain exact text representation
# (avoid Excel auto-converting numeric-looking strings to numbers / losing trailing zeros, double dots, etc.)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '43.753.46'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').Value = '2.289.71'
$ws.Range('E3').Value = '  -1.96%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '96.67'
$ws.Range('E5').Value = '  +2.27%  '
$ws.Range('D6').Value = '269.37'
$ws.Range('E6').Value = '  -0.33%  '
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  -1.38%  '
$ws.Range('D10').Value = '45.39'
$ws.Range('E10').Value = '  +1.55%  '
$ws.Range('E11').Value = '  -0.37%  '
$ws.Range('D12').Value = '7.94'
$ws.Range('E12').Value = '  -1.24%  '
$ws.Range('D14').Value = '15.70'
$ws.Range('E14').Value = '  +0.42%  '
$ws.Range('D15').Value = '2.632.80'
$ws.Range('E15').Value = '  -2.03%  '
$ws.Range('D16').Value = '0.857'
$ws.Range('E16').Value = '  +0.00%  '
$ws.Range('D17').Value = '2.296.60'
$ws.Range('E17').Value = '  -1.55%  '
$ws.Range('D18').Value = '43.752.70'
$ws.Range('E18').Value = '  +0.20%  '
$ws.Range('E19').Value = '  +3.83%  '
$ws.Range('D20').Value = '6.19'
$ws.Range('E20').Value = '  -2.82%  '
$ws.Range('D21').Value = '72.17'
$ws.Range('E21').Value = '  +0.56%  '
$ws.Range('D22').Value = '2.48'
$ws.Range('E22').Value = '  +10.05%  '
$ws.Range('D23').Value = '232.66'
$ws.Range('E23').Value = '  -2.30%  '
$ws.Range('E24').Value = '  -5.05%  '
$ws.Range('D25').Value = '2.71'
$ws.Range('E25').Value = '  +7.69%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '2.28'
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = '38.59'
$ws.Range('E30').Value = '  +0.27%  '
$ws.Range('D31').Value = '175.21'
$ws.Range('E31').Value = '  +1.91%  '
$ws.Range('D32').Value = '21.83'
$ws.Range('E32').Value = '  -3.62%  '
$ws.Range('D33').Value = '0.0897'
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('D34').Value = '5.44'
$ws.Range('E34').Value = '  -0.66%  '
$ws.Range('E35').Value = '  +0.35%  '
$ws.Range('D36').Value = '4.63'
$ws.Range('E36').Value = '  +6.41%  '
$ws.Range('E37').Value = '  +0.28%  '
$ws.Range('E38').Value = '  -1.09%  '
$ws.Range('D39').Value = '3.61'
$ws.Range('E39').Value = '  +6.56%  '
$ws.Range('E40').Value = '  +1.40%  '
$ws.Range('E41').Value = '  -2.46%  '
$ws.Range('D42').Value = '12.29'
$ws.Range('E42').Value = '  +2.25%  '
$ws.Range('E43').Value = '  -2.32%  '
$ws.Range('D44').Value = '64.66'
$ws.Range('E44').Value = '  +4.97%  '
$ws.Range('B45').Value = 'THORChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D45').Value = '5.20'
$ws.Range('E45').Value = '  -3.08%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '8.72'
$ws.Range('E46').Value = '  -3.98%  '
$ws.Range('E47').Value = '  -0.22%  '
$ws.Range('E48').Value = '  -0.52%  '
$ws.Range('D49').Value = '97.53'
$ws.Range('E49').Value = '  -2.88%  '
$ws.Range('B50').Value = 'WOONetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D50').Value = '0.438'
$ws.Range('E50').Value = '  +5.98%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value = '1.50'
$ws.Range('E51').Value = '  +9.89%  '
